# Update column G ("K", formerly "Strike#") values for rows 2-9 on Sheet1.
# Per commit message: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals" - the K column values are recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 2
